$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value = 12089.818
$ws.Cells.Item(18, 9).Value = 12554.223
$ws.Cells.Item(18, 10).Value = 10000
$ws.Cells.Item(18, 11).Value = 12554.223
$ws.Cells.Item(18, 12).Value = 10000
$ws.Cells.Item(18, 13).Value = -12270.223
$ws.Cells.Item(18, 14).Value = -10568
$ws.Cells.Item(70, 8).Value = 5033.9165
$ws.Cells.Item(70, 10).Value = 5120.5
$ws.Cells.Item(70, 12).Value = 15361.5
$ws.Cells.Item(70, 14).Value = -15901.5
$ws.Cells.Item(73, 8).Value = 5033.9165
$ws.Cells.Item(73, 10).Value = 5120.5
$ws.Cells.Item(73, 12).Value = 15361.5
$ws.Cells.Item(73, 14).Value = -17233.5
$ws.Cells.Item(116, 8).Value = 7229.8335
$ws.Cells.Item(116, 9).Value = 6481.4
$ws.Cells.Item(116, 11).Value = 6481.4
$ws.Cells.Item(116, 13).Value = -3039.4
$ws.Cells.Item(129, 8).Value = 2376.5908
$ws.Cells.Item(129, 9).Value = 1225.9231
$ws.Cells.Item(129, 10).Value = 4038.6667
$ws.Cells.Item(129, 11).Value = 3677.7693
$ws.Cells.Item(129, 12).Value = 12116.0001
$ws.Cells.Item(129, 13).Value = 1322.2307
$ws.Cells.Item(129, 14).Value = -22116.0001
$ws.Cells.Item(132, 8).Value = 2663.8
$ws.Cells.Item(132, 9).Value = 2655.6538
$ws.Cells.Item(132, 11).Value = 7966.9614
$ws.Cells.Item(132, 13).Value = -5436.9614
$ws.Cells.Item(141, 8).Value = 4492.467
$ws.Cells.Item(141, 9).Value = 4627.643
$ws.Cells.Item(141, 11).Value = 13882.929
$ws.Cells.Item(141, 13).Value = -8702.929

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(110, 8).Value = 2302.1875
$ws.Cells.Item(110, 9).Value = 1889.96
$ws.Cells.Item(110, 10).Value = 3774.4285
$ws.Cells.Item(110, 11).Value = 1889.96
$ws.Cells.Item(110, 12).Value = 3774.4285
$ws.Cells.Item(110, 13).Value = 155.04
$ws.Cells.Item(110, 14).Value = -7864.4285

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 1766.65
$ws.Cells.Item(20, 9).Value = 1697.4615
$ws.Cells.Item(20, 11).Value = 1697.4615
$ws.Cells.Item(20, 13).Value = -1450.4615
$ws.Cells.Item(82, 8).Value = 32697.842
$ws.Cells.Item(82, 9).Value = 3941.7778
$ws.Cells.Item(82, 10).Value = 58578.3
$ws.Cells.Item(82, 11).Value = 3941.7778
$ws.Cells.Item(82, 12).Value = 58578.3
$ws.Cells.Item(82, 13).Value = -3558.7778
$ws.Cells.Item(82, 14).Value = -59344.3
$ws.Cells.Item(85, 8).Value = 32697.842
$ws.Cells.Item(85, 9).Value = 3941.7778
$ws.Cells.Item(85, 10).Value = 58578.3
$ws.Cells.Item(85, 11).Value = 3941.7778
$ws.Cells.Item(85, 12).Value = 58578.3
$ws.Cells.Item(85, 13).Value = -2615.7778
$ws.Cells.Item(85, 14).Value = -61230.3
$ws.Cells.Item(86, 8).Value = 125015760
$ws.Cells.Item(86, 9).Value = 125015760
$ws.Cells.Item(86, 11).Value = 125015760
$ws.Cells.Item(86, 13).Value = -125014637
$ws.Cells.Item(89, 8).Value = 125015760
$ws.Cells.Item(89, 9).Value = 125015760
$ws.Cells.Item(89, 11).Value = 625078800
$ws.Cells.Item(89, 13).Value = -625073184
$ws.Cells.Item(99, 8).Value = 2470.1
$ws.Cells.Item(99, 9).Value = 2470.1
$ws.Cells.Item(99, 11).Value = 2470.1
$ws.Cells.Item(99, 13).Value = -972.0999999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(17, 8).Value = 15033.714
$ws.Cells.Item(17, 9).Value = 1004.8
$ws.Cells.Item(17, 11).Value = 1004.8
$ws.Cells.Item(17, 13).Value = -830.8
$ws.Cells.Item(41, 8).Value = 22255.084
$ws.Cells.Item(41, 9).Value = 17000.5
$ws.Cells.Item(41, 11).Value = 17000.5
$ws.Cells.Item(41, 13).Value = -16572.5
$ws.Cells.Item(50, 8).Value = 32725.857
$ws.Cells.Item(50, 9).Value = 26999
$ws.Cells.Item(50, 11).Value = 26999
$ws.Cells.Item(50, 13).Value = -26374
$ws.Cells.Item(51, 8).Value = 23653.2
$ws.Cells.Item(51, 9).Value = 21000
$ws.Cells.Item(51, 11).Value = 21000
$ws.Cells.Item(51, 13).Value = -20264
$ws.Cells.Item(53, 8).Value = 35444.332
$ws.Cells.Item(53, 10).Value = 35444.332
$ws.Cells.Item(53, 12).Value = 35444.332
$ws.Cells.Item(53, 14).Value = -36658.332
$ws.Cells.Item(60, 8).Value = 49518.766
$ws.Cells.Item(60, 9).Value = 29500
$ws.Cells.Item(60, 11).Value = 29500
$ws.Cells.Item(60, 13).Value = -28989
$ws.Cells.Item(61, 8).Value = 23653.2
$ws.Cells.Item(61, 9).Value = 21000
$ws.Cells.Item(61, 11).Value = 21000
$ws.Cells.Item(61, 13).Value = -20652
$ws.Cells.Item(68, 8).Value = 56995.4
$ws.Cells.Item(68, 10).Value = 56995.4
$ws.Cells.Item(68, 12).Value = 56995.4
$ws.Cells.Item(68, 14).Value = -58493.4
$ws.Cells.Item(71, 8).Value = 56995.4
$ws.Cells.Item(71, 10).Value = 56995.4
$ws.Cells.Item(71, 12).Value = 170986.2
$ws.Cells.Item(71, 14).Value = -178474.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(68, 8).Value = 1359
$ws.Cells.Item(68, 9).Value = 1045.75
$ws.Cells.Item(68, 11).Value = 3137.25
$ws.Cells.Item(68, 13).Value = -2326.25
$ws.Cells.Item(71, 8).Value = 1359
$ws.Cells.Item(71, 9).Value = 1045.75
$ws.Cells.Item(71, 11).Value = 9411.75
$ws.Cells.Item(71, 13).Value = -5355.75
$ws.Cells.Item(86, 8).Value = 1720.3636
$ws.Cells.Item(86, 9).Value = 1189.4
$ws.Cells.Item(86, 10).Value = 2162.8333
$ws.Cells.Item(86, 11).Value = 3568.2
$ws.Cells.Item(86, 12).Value = 6488.499899999999
$ws.Cells.Item(86, 13).Value = -2382.2
$ws.Cells.Item(86, 14).Value = -8860.499899999999
$ws.Cells.Item(89, 8).Value = 1720.3636
$ws.Cells.Item(89, 9).Value = 1189.4
$ws.Cells.Item(89, 10).Value = 2162.8333
$ws.Cells.Item(89, 11).Value = 10704.6
$ws.Cells.Item(89, 12).Value = 19465.4997
$ws.Cells.Item(89, 13).Value = -4776.6
$ws.Cells.Item(89, 14).Value = -31321.4997

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 5058.7744
$ws.Cells.Item(102, 9).Value = 4661.115
$ws.Cells.Item(102, 10).Value = 7126.6
$ws.Cells.Item(102, 11).Value = 4661.115
$ws.Cells.Item(102, 12).Value = 7126.6
$ws.Cells.Item(102, 13).Value = -3039.115
$ws.Cells.Item(102, 14).Value = -10370.6
$ws.Cells.Item(122, 8).Value = 40130.63
$ws.Cells.Item(122, 9).Value = 61887.176
$ws.Cells.Item(122, 10).Value = 3144.5
$ws.Cells.Item(122, 11).Value = 185661.528
$ws.Cells.Item(122, 12).Value = 9433.5
$ws.Cells.Item(122, 13).Value = -183211.528
$ws.Cells.Item(122, 14).Value = -14333.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1159.2354
$ws.Cells.Item(22, 9).Value = 575.4545000000001
$ws.Cells.Item(22, 10).Value = 2229.5
$ws.Cells.Item(22, 11).Value = 575.4545000000001
$ws.Cells.Item(22, 12).Value = 2229.5
$ws.Cells.Item(22, 13).Value = -280.4545000000001
$ws.Cells.Item(22, 14).Value = -2819.5
$ws.Cells.Item(27, 8).Value = 1159.2354
$ws.Cells.Item(27, 9).Value = 575.4545000000001
$ws.Cells.Item(27, 10).Value = 2229.5
$ws.Cells.Item(27, 11).Value = 575.4545000000001
$ws.Cells.Item(27, 12).Value = 2229.5
$ws.Cells.Item(27, 13).Value = -468.4545000000001
$ws.Cells.Item(27, 14).Value = -2443.5
$ws.Cells.Item(40, 8).Value = 15792.914
$ws.Cells.Item(40, 9).Value = 17601.896
$ws.Cells.Item(40, 11).Value = 17601.896
$ws.Cells.Item(40, 13).Value = -17465.896
$ws.Cells.Item(46, 8).Value = 2403.2
$ws.Cells.Item(46, 9).Value = 628.5
$ws.Cells.Item(46, 11).Value = 628.5
$ws.Cells.Item(46, 13).Value = -440.5
$ws.Cells.Item(48, 8).Value = 35045.668
$ws.Cells.Item(48, 10).Value = 35045.668
$ws.Cells.Item(48, 12).Value = 35045.668
$ws.Cells.Item(48, 14).Value = -36367.668
$ws.Cells.Item(122, 8).Value = 3514.3667
$ws.Cells.Item(122, 9).Value = 3339.5
$ws.Cells.Item(122, 10).Value = 4651
$ws.Cells.Item(122, 11).Value = 10018.5
$ws.Cells.Item(122, 12).Value = 13953
$ws.Cells.Item(122, 13).Value = -7568.5
$ws.Cells.Item(122, 14).Value = -18853
$ws.Cells.Item(132, 8).Value = 2168822
$ws.Cells.Item(132, 9).Value = 2529117.5
$ws.Cells.Item(132, 10).Value = 7049.6665
$ws.Cells.Item(132, 11).Value = 7587352.5
$ws.Cells.Item(132, 12).Value = 21148.9995
$ws.Cells.Item(132, 13).Value = -7584822.5
$ws.Cells.Item(132, 14).Value = -26208.9995
$ws.Cells.Item(136, 8).Value = 9262598
$ws.Cells.Item(136, 9).Value = 11497389
$ws.Cells.Item(136, 10).Value = 4178.143
$ws.Cells.Item(136, 11).Value = 34492167
$ws.Cells.Item(136, 12).Value = 12534.429
$ws.Cells.Item(136, 13).Value = -34489617
$ws.Cells.Item(136, 14).Value = -17634.429

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 1997.0769
$ws.Cells.Item(96, 9).Value = 2490
$ws.Cells.Item(96, 10).Value = 1689
$ws.Cells.Item(96, 11).Value = 2490
$ws.Cells.Item(96, 12).Value = 1689
$ws.Cells.Item(96, 13).Value = -1117
$ws.Cells.Item(96, 14).Value = -4435
